# İş Takip Güncellemesi - 05.08.2025 14:41:39
# Update the "Güncelleme" sheet: the "UÇUŞ TARİHİ(YAPILAN)" (flight
# date) is now recorded for several units, so the date that had been
# entered under "ARAZİ YERSEL ÖLÇÜM TARİHİ(YAPILAN)" is moved to the
# flight-date column, the flight-evaluation / evaluation status columns
# are cleared pending the new evaluation, and postings that have been
# finalized are marked "Yapıldı".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Güncelleme")

# Row 2
$ws.Range("P2").Value = "Yapıldı"

# Row 3
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = "'2024-11-05"
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = ""

# Row 4
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = "'2024-11-05"
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("N4").Value = ""
$ws.Range("P4").Value = "Yapıldı"

# Row 6
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = "'2024-11-07"
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = ""

# Row 8
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = "'2024-11-07"
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = ""
$ws.Range("N8").Value = ""
$ws.Range("P8").Value = "Yapıldı"

# Row 10
$ws.Range("I10").Value = ""
$ws.Range("J10").Value = "'2024-11-07"
$ws.Range("K10").Value = ""
$ws.Range("L10").Value = ""

# Row 23
$ws.Range("I23").Value = ""
$ws.Range("J23").Value = "'2024-11-11"
